$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last 3 rows (44-46) since the new dataset has fewer rows
$ws.Range("A44:C46").Delete()

# Update rows 2-43 with the new dataset (shifted to 2025-10-17 / serial 45947)
$ws.Cells.Item(2, 1).Value = 45947
$ws.Cells.Item(2, 2).Value = 0.098
$ws.Cells.Item(2, 3).Value = 11.809
$ws.Cells.Item(3, 1).Value = 45947.01041666666
$ws.Cells.Item(3, 2).Value = 0
$ws.Cells.Item(3, 3).Value = 10.638
$ws.Cells.Item(4, 1).Value = 45947.02083333334
$ws.Cells.Item(4, 2).Value = 0
$ws.Cells.Item(4, 3).Value = 7.913
$ws.Cells.Item(5, 1).Value = 45947.03125
$ws.Cells.Item(5, 2).Value = 0
$ws.Cells.Item(5, 3).Value = 8.046
$ws.Cells.Item(6, 1).Value = 45947.04166666666
$ws.Cells.Item(6, 2).Value = 0
$ws.Cells.Item(6, 3).Value = 19.246
$ws.Cells.Item(7, 1).Value = 45947.05208333334
$ws.Cells.Item(7, 2).Value = 0.176
$ws.Cells.Item(7, 3).Value = 7.18
$ws.Cells.Item(8, 1).Value = 45947.0625
$ws.Cells.Item(8, 2).Value = 0
$ws.Cells.Item(8, 3).Value = 11.864
$ws.Cells.Item(9, 1).Value = 45947.07291666666
$ws.Cells.Item(9, 2).Value = 0.065
$ws.Cells.Item(9, 3).Value = 8.075
$ws.Cells.Item(10, 1).Value = 45947.08333333334
$ws.Cells.Item(10, 2).Value = 0.582
$ws.Cells.Item(10, 3).Value = 7.819
$ws.Cells.Item(11, 1).Value = 45947.09375
$ws.Cells.Item(11, 2).Value = 2.663
$ws.Cells.Item(11, 3).Value = 0.375
$ws.Cells.Item(12, 1).Value = 45947.10416666666
$ws.Cells.Item(12, 2).Value = 2.323
$ws.Cells.Item(12, 3).Value = 3.345
$ws.Cells.Item(13, 1).Value = 45947.11458333334
$ws.Cells.Item(13, 2).Value = 0.912
$ws.Cells.Item(13, 3).Value = 3.041
$ws.Cells.Item(14, 1).Value = 45947.125
$ws.Cells.Item(14, 2).Value = 3.634
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(15, 1).Value = 45947.13541666666
$ws.Cells.Item(15, 2).Value = 4.996
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(16, 1).Value = 45947.14583333334
$ws.Cells.Item(16, 2).Value = 1.748
$ws.Cells.Item(16, 3).Value = 1.14
$ws.Cells.Item(17, 1).Value = 45947.15625
$ws.Cells.Item(17, 2).Value = 1.736
$ws.Cells.Item(17, 3).Value = 0.135
$ws.Cells.Item(18, 1).Value = 45947.16666666666
$ws.Cells.Item(18, 2).Value = 0.002
$ws.Cells.Item(18, 3).Value = 10.328
$ws.Cells.Item(19, 1).Value = 45947.17708333334
$ws.Cells.Item(19, 2).Value = 2.952
$ws.Cells.Item(19, 3).Value = 0.152
$ws.Cells.Item(20, 1).Value = 45947.1875
$ws.Cells.Item(20, 2).Value = 0.801
$ws.Cells.Item(20, 3).Value = 2.126
$ws.Cells.Item(21, 1).Value = 45947.19791666666
$ws.Cells.Item(21, 2).Value = 1.494
$ws.Cells.Item(21, 3).Value = 0.033
$ws.Cells.Item(22, 1).Value = 45947.20833333334
$ws.Cells.Item(22, 2).Value = 1.973
$ws.Cells.Item(22, 3).Value = 0.608
$ws.Cells.Item(23, 1).Value = 45947.21875
$ws.Cells.Item(23, 2).Value = 7.344
$ws.Cells.Item(23, 3).Value = 0.062
$ws.Cells.Item(24, 1).Value = 45947.22916666666
$ws.Cells.Item(24, 2).Value = 3.491
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(25, 1).Value = 45947.23958333334
$ws.Cells.Item(25, 2).Value = 4.25
$ws.Cells.Item(25, 3).Value = 5.422
$ws.Cells.Item(26, 1).Value = 45947.25
$ws.Cells.Item(26, 2).Value = 5.118
$ws.Cells.Item(26, 3).Value = 0
$ws.Cells.Item(27, 1).Value = 45947.26041666666
$ws.Cells.Item(27, 2).Value = 0.843
$ws.Cells.Item(27, 3).Value = 0
$ws.Cells.Item(28, 1).Value = 45947.27083333334
$ws.Cells.Item(28, 2).Value = 2.549
$ws.Cells.Item(28, 3).Value = 0
$ws.Cells.Item(29, 1).Value = 45947.28125
$ws.Cells.Item(29, 2).Value = 3.773
$ws.Cells.Item(29, 3).Value = 0
$ws.Cells.Item(30, 1).Value = 45947.29166666666
$ws.Cells.Item(30, 2).Value = 6.059
$ws.Cells.Item(30, 3).Value = 0
$ws.Cells.Item(31, 1).Value = 45947.30208333334
$ws.Cells.Item(31, 2).Value = 8.823
$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(32, 1).Value = 45947.3125
$ws.Cells.Item(32, 2).Value = 4.558
$ws.Cells.Item(32, 3).Value = 13.015
$ws.Cells.Item(33, 1).Value = 45947.32291666666
$ws.Cells.Item(33, 2).Value = 0
$ws.Cells.Item(33, 3).Value = 88.074
$ws.Cells.Item(34, 1).Value = 45947.33333333334
$ws.Cells.Item(34, 2).Value = 0
$ws.Cells.Item(34, 3).Value = 29.964
$ws.Cells.Item(35, 1).Value = 45947.34375
$ws.Cells.Item(35, 2).Value = 0
$ws.Cells.Item(35, 3).Value = 33.722
$ws.Cells.Item(36, 1).Value = 45947.35416666666
$ws.Cells.Item(36, 2).Value = 0
$ws.Cells.Item(36, 3).Value = 31.844
$ws.Cells.Item(37, 1).Value = 45947.36458333334
$ws.Cells.Item(37, 2).Value = 0
$ws.Cells.Item(37, 3).Value = 29.284
$ws.Cells.Item(38, 1).Value = 45947.375
$ws.Cells.Item(38, 2).Value = 0.003
$ws.Cells.Item(38, 3).Value = 19.905
$ws.Cells.Item(39, 1).Value = 45947.38541666666
$ws.Cells.Item(39, 2).Value = 0
$ws.Cells.Item(39, 3).Value = 40.912
$ws.Cells.Item(40, 1).Value = 45947.39583333334
$ws.Cells.Item(40, 2).Value = 2.186
$ws.Cells.Item(40, 3).Value = 4.816
$ws.Cells.Item(41, 1).Value = 45947.40625
$ws.Cells.Item(41, 2).Value = 0.403
$ws.Cells.Item(41, 3).Value = 7.674
$ws.Cells.Item(42, 1).Value = 45947.41666666666
$ws.Cells.Item(42, 2).Value = 7.12
$ws.Cells.Item(42, 3).Value = 0.45
$ws.Cells.Item(43, 1).Value = 45947.42708333334
$ws.Cells.Item(43, 2).Value = 4.468
$ws.Cells.Item(43, 3).Value = 10.326

Write-Host "Applied IGCC Netting Flows update"
